# Daily attendance processing - sort "Recorded By" (column G) entries alphabetically
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $val -ne "") {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $sorted = $trimmed | Sort-Object { $_.ToLower() }
        $newVal = [string]::Join(", ", $sorted)

        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
